$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle 2")

$ws.Range("C6").Value = "schwarz"
$ws.Range("C7").Value = "weiß"
$ws.Range("C8").Value = "grün"
$ws.Range("C9").Value = "blau"
$ws.Range("C10").Value = "gelb"
$ws.Range("C11").Value = "blau"
$ws.Range("C12").Value = "weiß"
$ws.Range("C17").Value = "gelb"
$ws.Range("C18").Value = "schwarz"
